$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting the existing row 6
# (and anything below it) down by one row.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new price record.
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44551
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100101
$ws.Cells.Item(6, 8).Value = "Berries"
$ws.Cells.Item(6, 9).Value = 100101004
$ws.Cells.Item(6, 10).Value = "Frambuesa"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 120
$ws.Cells.Item(6, 14).Value = 4500
$ws.Cells.Item(6, 15).Value = 4500
$ws.Cells.Item(6, 16).Value = 4500
$ws.Cells.Item(6, 17).Value = "$/envase 1 kilo"
$ws.Cells.Item(6, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(6, 19).Value = 4500
$ws.Cells.Item(6, 20).Value = 1

$wb.Save()
